$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. The sheet currently named "总计" (last sheet, sheetId 6) becomes the
#    new "2022-Q1" detail sheet. Before we touch it, clone its current
#    ("总计" summary) content into a brand new trailing sheet that will
#    be renamed back to "总计" - this preserves the sheetId/order
#    progression (…,5,6,7) the diff expects.
# ---------------------------------------------------------------------
$oldTotal = $wb.Worksheets.Item("总计")

$newTotal = $wb.Worksheets.Add($null, $oldTotal)
$oldTotal.UsedRange.Copy($newTotal.Range("A1"))
$oldTotal.Name = "2022-Q1"
$newTotal.Name = "总计"

# Insert a fresh row for the new "2022-Q1" summary entry right under the
# header row, push the rest (2021-Q4 .. 2020-Q4) down by one.
$newTotal.Rows(2).Insert()
$newTotal.Range("A2:D2").ClearFormats()
$newTotal.Range("A3").Copy()
$newTotal.Range("A2").PasteSpecial(-4122)

$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 2
$newTotal.Range("D2").Value = 1.31

# Re-number the pandas-style index column (A) 0..5 top to bottom.
$newTotal.Range("A2").Value = 0
$newTotal.Range("A3").Value = 1
$newTotal.Range("A4").Value = 2
$newTotal.Range("A5").Value = 3
$newTotal.Range("A6").Value = 4
$newTotal.Range("A7").Value = 5

# ---------------------------------------------------------------------
# 2. Turn the original "总计" sheet into the "2022-Q1" fund-holdings
#    detail sheet (same column layout as the other quarterly sheets).
# ---------------------------------------------------------------------
$oldTotal.Cells.Clear()
$template = $wb.Worksheets.Item("2021-Q4")
$template.UsedRange.Copy($oldTotal.Range("A1"))
$oldTotal.Rows(4).Delete()

$oldTotal.Range("B2").NumberFormat = "@"
$oldTotal.Range("B2").Value = "000727"
$oldTotal.Range("C2").Value = "融通健康产业灵活配置混合A"
$oldTotal.Range("D2").NumberFormat = "@"
$oldTotal.Range("D2").Value = "15.30"
$oldTotal.Range("E2").NumberFormat = "@"
$oldTotal.Range("E2").Value = "94.68"
$oldTotal.Range("F2").NumberFormat = "@"
$oldTotal.Range("F2").Value = "7.08"
$oldTotal.Range("G2").NumberFormat = "@"
$oldTotal.Range("G2").Value = "1.0832"
$oldTotal.Range("H2").Value = 5

$oldTotal.Range("B3").NumberFormat = "@"
$oldTotal.Range("B3").Value = "009274"
$oldTotal.Range("C3").Value = "融通健康产业灵活配置混合C"
$oldTotal.Range("D3").NumberFormat = "@"
$oldTotal.Range("D3").Value = "3.16"
$oldTotal.Range("E3").NumberFormat = "@"
$oldTotal.Range("E3").Value = "94.68"
$oldTotal.Range("F3").NumberFormat = "@"
$oldTotal.Range("F3").Value = "7.08"
$oldTotal.Range("G3").NumberFormat = "@"
$oldTotal.Range("G3").Value = "0.2237"
$oldTotal.Range("H3").Value = 5
